$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (Question/Answer shift right by one column)
$ws.Range("B:B").Insert()

# New header for the inserted "S No" column
$ws.Range("B2").Value = "S No"
$ws.Range("B2").Style = $ws.Range("C2").Style

# Fill serial numbers 1..10 for the data rows (rows 3-12)
for ($i = 0; $i -lt 10; $i++) {
    $row = 3 + $i
    $cell = $ws.Cells.Item($row, 2)
    $cell.Value = $i + 1
    $cell.Style = $ws.Cells.Item($row, 3).Style
}

# Resize the AutoFilter to cover the new column range
$ws.Range("B2:D12").AutoFilter(1) | Out-Null

# Adjust selection similar to the saved workbook state
$ws.Range("C7").Select()

$wb.Save()
